# Updated symbol list on Tue Jan 17 05:22:10 UTC 2023 with GitHub Actions
# Refresh the "Price" (D) and "Volume(1h)" (E) columns on the crypto sheet.
# Values are stored as literal text in the workbook (e.g. "299.28", "-1.82%"),
# so each write is forced to text via a leading apostrophe and the cell style
# is reset to "Normal" afterwards so no stray number-format/quote-prefix style
# sticks to the cell (keeps formatting identical to the source).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2";  Value = "298.93" },
    @{ Cell = "E2";  Value = "-1.98%" },
    @{ Cell = "D3";  Value = "31.33" },
    @{ Cell = "E3";  Value = "-1.48%" },
    @{ Cell = "D4";  Value = "5.104" },
    @{ Cell = "E4";  Value = "-1.38%" },
    @{ Cell = "D5";  Value = "0.07961" },
    @{ Cell = "E5";  Value = "6.05%" },
    @{ Cell = "D6";  Value = "2.281" },
    @{ Cell = "E6";  Value = "-1.52%" },
    @{ Cell = "D7";  Value = "7.773" },
    @{ Cell = "E7";  Value = "-3.10%" },
    @{ Cell = "D8";  Value = "3.861" },
    @{ Cell = "E8";  Value = "-0.45%" },
    @{ Cell = "D9";  Value = "0.9231" },
    @{ Cell = "E9";  Value = "0.79%" },
    @{ Cell = "D10"; Value = "0.1735" },
    @{ Cell = "E10"; Value = "0.12%" },
    @{ Cell = "D11"; Value = "0.07539" },
    @{ Cell = "E11"; Value = "-2.04%" },
    @{ Cell = "D12"; Value = "0.09492" },
    @{ Cell = "E12"; Value = "15.42%" },
    @{ Cell = "D13"; Value = "0.03024" },
    @{ Cell = "E13"; Value = "-0.60%" },
    @{ Cell = "D14"; Value = "0.1004" },
    @{ Cell = "E14"; Value = "0.87%" },
    @{ Cell = "D15"; Value = "0.001508" },
    @{ Cell = "E15"; Value = "-0.35%" },
    @{ Cell = "D16"; Value = "0.005862" },
    @{ Cell = "E16"; Value = "-4.26%" },
    @{ Cell = "D17"; Value = "3.488" },
    @{ Cell = "E17"; Value = "-0.41%" },
    @{ Cell = "D18"; Value = "2.268" },
    @{ Cell = "E18"; Value = "1.15%" },
    @{ Cell = "E20"; Value = "0.47%" },
    @{ Cell = "D21"; Value = "3.904" },
    @{ Cell = "E21"; Value = "-16.15%" },
    @{ Cell = "D22"; Value = "0.1699" },
    @{ Cell = "E22"; Value = "8.62%" },
    @{ Cell = "E23"; Value = "-0.20%" },
    @{ Cell = "D24"; Value = "0.001247" },
    @{ Cell = "E24"; Value = "-1.23%" },
    @{ Cell = "D25"; Value = "0.004489" },
    @{ Cell = "D26"; Value = "0.0001200" },
    @{ Cell = "E26"; Value = "-7.62%" },
    @{ Cell = "D27"; Value = "0.0003394" },
    @{ Cell = "E27"; Value = "24.05%" },
    @{ Cell = "D39"; Value = "0.01763" },
    @{ Cell = "E39"; Value = "0.21%" },
    @{ Cell = "E40"; Value = "1.07%" },
    @{ Cell = "D41"; Value = "0.006966" },
    @{ Cell = "E41"; Value = "-4.89%" },
    @{ Cell = "D42"; Value = "0.1360" },
    @{ Cell = "E42"; Value = "-0.30%" },
    @{ Cell = "D43"; Value = "0.002189" },
    @{ Cell = "E43"; Value = "2.42%" },
    @{ Cell = "D44"; Value = "0.01015" },
    @{ Cell = "E44"; Value = "-6.08%" },
    @{ Cell = "D45"; Value = "0.00006312" },
    @{ Cell = "E45"; Value = "-2.75%" },
    @{ Cell = "E46"; Value = "-0.13%" },
    @{ Cell = "D47"; Value = "0.007971" },
    @{ Cell = "E47"; Value = "-19.29%" },
    @{ Cell = "E48"; Value = "41.01%" },
    @{ Cell = "D49"; Value = "0.00002097" },
    @{ Cell = "E49"; Value = "-0.13%" },
    @{ Cell = "E50"; Value = "-0.06%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Leading apostrophe forces Excel to store this as literal text rather
    # than auto-converting "298.93" / "-1.98%" into a number/percentage.
    $cell.Value = "'" + $u.Value
    # Reset style so the quote-prefix formatting doesn't linger on the cell.
    $cell.Style = "Normal"
}
